$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 361, pushing all rows
# from 361-375 down to 363-377.
$ws.Rows("361:362").Insert()

# New row 361: Conconina(o), Segunda, week of 2021-11-09
$ws.Range("A361").Value = 4
$ws.Range("B361").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C361").Value = "Los Lagos"
$ws.Range("D361").Value = 44509
$ws.Range("E361").Value = 10
$ws.Range("F361").Value = 100112033
$ws.Range("G361").Value = "Lechuga"
$ws.Range("H361").Value = "Conconina(o)"
$ws.Range("I361").Value = "Segunda"
$ws.Range("J361").Value = 120
$ws.Range("K361").Value = 8000
$ws.Range("L361").Value = 8000
$ws.Range("M361").Value = 8000
$ws.Range("N361").Value = "`$/caja 12 unidades"
$ws.Range("O361").Value = "Región Metropolitana"
$ws.Range("P361").Value = 667
$ws.Range("Q361").Value = 12
$ws.Range("R361").Value = "Hortaliza"

# New row 362: Escarola, Primera, week of 2021-11-09
$ws.Range("A362").Value = 4
$ws.Range("B362").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C362").Value = "Los Lagos"
$ws.Range("D362").Value = 44509
$ws.Range("E362").Value = 10
$ws.Range("F362").Value = 100112033
$ws.Range("G362").Value = "Lechuga"
$ws.Range("H362").Value = "Escarola"
$ws.Range("I362").Value = "Primera"
$ws.Range("J362").Value = 400
$ws.Range("K362").Value = 8500
$ws.Range("L362").Value = 9000
$ws.Range("M362").Value = 8750
$ws.Range("N362").Value = "`$/caja 15 unidades"
$ws.Range("O362").Value = "Región de Coquimbo"
$ws.Range("P362").Value = 583
$ws.Range("Q362").Value = 15
$ws.Range("R362").Value = "Hortaliza"
